# Lesson-05 update:
#  1. Header/footer "date" placeholder text 11/6/2012 -> 11/8/2012
#     (slide master + every custom layout that carries the placeholder).
#  2. Slide 14 ("Practice") task-description textbox rewritten: the long
#     in-slide task text is replaced by a short pointer to the new
#     standalone task file (task-gui-seconds-counter.docx).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Fixed "date" header/footer field: 11/6/2012 -> 11/8/2012
# ---------------------------------------------------------------------
$newDate = "11/8/2012"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------
# 2. Slide 14: replace the task description with a pointer to the
#    separate task file.
# ---------------------------------------------------------------------
$quoteOpen  = [char]0x201C   # "
$quoteClose = [char]0x201D   # "

$newText = "Создание класса " + $quoteOpen + "Секундомер" + $quoteClose + " - task-gui-seconds-counter.docx"

$slide14 = $p.Slides.Item(14)
for ($i = 1; $i -le $slide14.Shapes.Count; $i++) {
    $sh = $slide14.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 6") {

        $tr = $sh.TextFrame.TextRange
        $tr.Text = $newText

        # Re-create the original run boundaries (same text, same
        # formatting) so the paragraph keeps the same granular run
        # layout as before the edit:
        #   "Создание " | "класса " | """ | "Секундомер" | """ | " - " | "task-gui-seconds-counter.docx"
        $tr.Characters(10, 7).Font.Italic  = 1   # "класса "
        $tr.Characters(17, 1).Font.Italic  = 1   # opening curly quote
        $tr.Characters(18, 10).Font.Italic = 1   # "Секундомер"
        $tr.Characters(28, 1).Font.Italic  = 1   # closing curly quote
        $tr.Characters(29, 3).Font.Italic  = 1   # " - "
        $tr.Characters(32, 29).Font.Italic = 1   # "task-gui-seconds-counter.docx"

        break
    }
}
